# Atualizado por script em 11-11-2023 08:45
#
# This script:
#  1) Re-orders several match rows (10/11, 26/27/28, 44/45, 80/81) so that
#     the F:V "match" columns end up in the sequence the scraper produced on
#     its latest run (the A:E "index/metadata" columns stay put - only the
#     match payload moves between the existing row slots).
#  2) Appends three freshly scraped matches as new rows 91, 92, 93.
#
# NOTE: this COM-interop host only binds PowerShell function arguments
# positionally - named parameters (-Row 10 / -Row:10) silently come through
# empty - so every helper below is called with plain positional arguments.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-MatchRow {
    param(
        $Row,
        $Home,
        $HomeGoals,
        $Away,
        $AwayGoals,
        $HomeOpenOdds,
        $HomeOpenTime,
        $HomeCloseOdds,
        $HomeCloseTime,
        $DrawOpenOdds,
        $DrawOpenTime,
        $DrawCloseOdds,
        $DrawCloseTime,
        $AwayOpenOdds,
        $AwayOpenTime,
        $AwayCloseOdds,
        $AwayCloseTime,
        $Url
    )

    $ws.Cells.Item($Row, 6).Value = $Home          # F home
    $ws.Cells.Item($Row, 7).Value = $HomeGoals      # G home_ft_gols
    $ws.Cells.Item($Row, 8).Value = $Away           # H away
    $ws.Cells.Item($Row, 9).Value = $AwayGoals      # I away_ft_gols
    $ws.Cells.Item($Row, 10).Value = $HomeOpenOdds  # J home_opening_odds
    $ws.Cells.Item($Row, 11).Value = $HomeOpenTime  # K home_opening_data_hora
    $ws.Cells.Item($Row, 12).Value = $HomeCloseOdds # L home_closing_odds
    $ws.Cells.Item($Row, 13).Value = $HomeCloseTime # M home_closing_data_hora
    $ws.Cells.Item($Row, 14).Value = $DrawOpenOdds  # N draw_opening_odds
    $ws.Cells.Item($Row, 15).Value = $DrawOpenTime  # O draw_opening_data_hora
    $ws.Cells.Item($Row, 16).Value = $DrawCloseOdds # P draw_closing_odds
    $ws.Cells.Item($Row, 17).Value = $DrawCloseTime # Q draw_closing_data_hora
    $ws.Cells.Item($Row, 18).Value = $AwayOpenOdds  # R away_opening_odds
    $ws.Cells.Item($Row, 19).Value = $AwayOpenTime  # S away_opening_data_hora
    $ws.Cells.Item($Row, 20).Value = $AwayCloseOdds # T away_closing_odds
    $ws.Cells.Item($Row, 21).Value = $AwayCloseTime # U away_closing_data_hora
    $ws.Cells.Item($Row, 22).Value = $Url           # V url_partida
}

function Set-FullRow {
    param(
        $Row,
        $Indice,
        $DataPartida,
        $Home,
        $HomeGoals,
        $Away,
        $AwayGoals,
        $HomeOpenOdds,
        $HomeOpenTime,
        $HomeCloseOdds,
        $HomeCloseTime,
        $DrawOpenOdds,
        $DrawOpenTime,
        $DrawCloseOdds,
        $DrawCloseTime,
        $AwayOpenOdds,
        $AwayOpenTime,
        $AwayCloseOdds,
        $AwayCloseTime,
        $Url
    )

    $ws.Cells.Item($Row, 1).Value = $Indice   # A Indice
    $ws.Cells.Item($Row, 1).Font.Bold = $true
    $ws.Cells.Item($Row, 1).HorizontalAlignment = -4108  # xlCenter
    $ws.Cells.Item($Row, 1).VerticalAlignment = -4160    # xlTop
    $ws.Cells.Item($Row, 1).Borders.LineStyle = 1        # xlContinuous

    $ws.Cells.Item($Row, 2).Value = "denmark"        # B pais
    $ws.Cells.Item($Row, 3).Value = "1st-division"    # C torneio
    $ws.Cells.Item($Row, 4).Value = "2023-2024"       # D temporada

    $ws.Cells.Item($Row, 5).Value = $DataPartida      # E data_partida
    $ws.Cells.Item($Row, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"

    Set-MatchRow $Row $Home $HomeGoals $Away $AwayGoals `
        $HomeOpenOdds $HomeOpenTime $HomeCloseOdds $HomeCloseTime `
        $DrawOpenOdds $DrawOpenTime $DrawCloseOdds $DrawCloseTime `
        $AwayOpenOdds $AwayOpenTime $AwayCloseOdds $AwayCloseTime `
        $Url
}

# ---------------------------------------------------------------------------
# 1) Rows 10 <-> 11 swap
# ---------------------------------------------------------------------------
Set-MatchRow 10 "Kolding IF" 2 "Hillerod" 1 `
    1.9 "24/07/2023 03:42" 1.53 "29/07/2023 13:58" `
    3.59 "24/07/2023 03:42" 4.41 "29/07/2023 13:58" `
    3.99 "24/07/2023 03:42" 5.96 "29/07/2023 13:58" `
    "https://www.betexplorer.com/football/denmark/1st-division/kolding-if-hillerod/vo6xggQ5/"

Set-MatchRow 11 "Aalborg" 1 "Horsens" 0 `
    1.75 "24/07/2023 03:42" 1.6 "29/07/2023 13:59" `
    3.88 "24/07/2023 03:42" 4.39 "29/07/2023 13:59" `
    4.41 "24/07/2023 03:42" 5.04 "29/07/2023 13:59" `
    "https://www.betexplorer.com/football/denmark/1st-division/aalborg-horsens/pO6YgZAa/"

# ---------------------------------------------------------------------------
# 2) Rows 26 / 27 / 28 rotate (26<-27, 27<-28, 28<-26)
# ---------------------------------------------------------------------------
Set-MatchRow 26 "Horsens" 3 "Helsingor" 1 `
    2.01 "13/08/2023 22:12" 2.31 "18/08/2023 18:51" `
    3.79 "13/08/2023 22:12" 3.76 "18/08/2023 18:51" `
    3.2 "13/08/2023 22:12" 2.89 "18/08/2023 18:51" `
    "https://www.betexplorer.com/football/denmark/1st-division/horsens-helsingor/nTtj8PV9/"

Set-MatchRow 27 "B.93" 0 "Fredericia" 5 `
    2.04 "14/08/2023 04:12" 2.44 "18/08/2023 18:47" `
    3.66 "14/08/2023 04:12" 3.73 "18/08/2023 18:47" `
    3.43 "14/08/2023 04:12" 2.73 "18/08/2023 18:37" `
    "https://www.betexplorer.com/football/denmark/1st-division/boldklubben-1893-fredericia/8xrn95G3/"

Set-MatchRow 28 "Kolding IF" 1 "Sonderjyske" 3 `
    2.08 "14/08/2023 04:12" 2.99 "18/08/2023 18:30" `
    3.6 "14/08/2023 04:12" 3.45 "18/08/2023 18:31" `
    3.38 "14/08/2023 04:12" 2.38 "18/08/2023 18:30" `
    "https://www.betexplorer.com/football/denmark/1st-division/kolding-if-sonderjyske/UiWf7qoG/"

# ---------------------------------------------------------------------------
# 3) Rows 44 <-> 45 swap
# ---------------------------------------------------------------------------
Set-MatchRow 44 "Horsens" 0 "B.93" 0 `
    1.79 "28/08/2023 01:42" 1.59 "01/09/2023 18:57" `
    4.02 "28/08/2023 01:42" 4.54 "01/09/2023 18:57" `
    3.72 "28/08/2023 01:42" 4.97 "01/09/2023 18:57" `
    "https://www.betexplorer.com/football/denmark/1st-division/horsens-boldklubben-1893/MqUfOyM7/"

Set-MatchRow 45 "Hillerod" 2 "Sonderjyske" 2 `
    3.9 "28/08/2023 18:42" 4.16 "01/09/2023 18:58" `
    3.8 "28/08/2023 18:42" 3.85 "01/09/2023 18:58" `
    1.79 "28/08/2023 18:42" 1.83 "01/09/2023 18:58" `
    "https://www.betexplorer.com/football/denmark/1st-division/hillerod-sonderjyske/EwHtnuEE/"

# ---------------------------------------------------------------------------
# 4) Rows 80 <-> 81 swap
# ---------------------------------------------------------------------------
Set-MatchRow 80 "B.93" 0 "Sonderjyske" 4 `
    5.01 "22/10/2023 15:12" 8.03 "27/10/2023 18:58" `
    4.47 "22/10/2023 15:12" 5.46 "27/10/2023 18:58" `
    1.57 "22/10/2023 15:12" 1.34 "27/10/2023 18:50" `
    "https://www.betexplorer.com/football/denmark/1st-division/boldklubben-1893-sonderjyske/2PEpqWy7/"

Set-MatchRow 81 "Hobro" 2 "Koge" 1 `
    1.74 "22/10/2023 16:12" 1.65 "27/10/2023 18:51" `
    3.92 "22/10/2023 16:12" 4.26 "27/10/2023 18:51" `
    4.42 "22/10/2023 16:12" 4.85 "27/10/2023 18:51" `
    "https://www.betexplorer.com/football/denmark/1st-division/hobro-koge/UBGxoh7f/"

# ---------------------------------------------------------------------------
# 5) New rows 91, 92, 93 appended at the end
# ---------------------------------------------------------------------------
Set-FullRow 91 90 45240.79166666666 `
    "Horsens" 1 "Fredericia" 1 `
    2.87 "05/11/2023 14:12" 2.63 "10/11/2023 18:51" `
    3.67 "05/11/2023 14:12" 3.73 "10/11/2023 18:51" `
    2.22 "05/11/2023 14:12" 2.53 "10/11/2023 18:51" `
    "https://www.betexplorer.com/football/denmark/1st-division/horsens-fredericia/v1Az2prn/"

Set-FullRow 92 91 45240.79166666666 `
    "B.93" 2 "Helsingor" 0 `
    2.51 "05/11/2023 15:13" 2.38 "10/11/2023 18:51" `
    3.55 "05/11/2023 15:13" 3.59 "10/11/2023 18:51" `
    2.68 "05/11/2023 15:13" 2.89 "10/11/2023 18:51" `
    "https://www.betexplorer.com/football/denmark/1st-division/boldklubben-1893-helsingor/Qc6W2QSu/"

Set-FullRow 93 92 45240.79166666666 `
    "Koge" 1 "Sonderjyske" 4 `
    6.52 "03/11/2023 19:13" 6.31 "10/11/2023 18:56" `
    5 "03/11/2023 19:13" 5.04 "10/11/2023 18:56" `
    1.37 "03/11/2023 19:13" 1.44 "10/11/2023 18:56" `
    "https://www.betexplorer.com/football/denmark/1st-division/koge-sonderjyske/AoWMrNsO/"
